$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) - Statistical Institution block (row 11)
$ws.Range("B11").Value = "'66.02"
$ws.Range("C11").Value = "'3.13"
$ws.Range("D11").Value = "'69.15"

# SME Associations block (rows 33-36)
$ws.Range("B33").Value = "'52.29"
$ws.Range("C33").Value = "'2.81"

$ws.Range("C34").Value = "'35.56"
$ws.Range("D34").Value = "'70.96"

$ws.Range("B36").Value = "'94.75"
$ws.Range("D36").Value = "'99.86"

# Value added to the economy (% of total) block (row 40)
$ws.Range("B40").Value = "'18.92"
$ws.Range("C40").Value = "'36.09"
$ws.Range("D40").Value = "'55.02"
